$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B (roll_no) entirely - this shifts subject1..subject5 (and all
# row data beneath them) one column to the left (C->B, D->C, E->D, F->E, G->F)
$ws.Columns.Item(2).Delete()
